$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin' -> 'Bitcoin'
$ws.Cells.Item(2, 4).Value = '29.709.87'
$ws.Cells.Item(2, 5).Value = '  +1.79%  '

# Row 3: 'Ethereum' -> 'Ethereum'
$ws.Cells.Item(3, 4).Value = '1.853.48'
$ws.Cells.Item(3, 5).Value = '  +1.37%  '

# Row 4: 'TetherUSD' -> 'TetherUSD'
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.9997'
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

# Row 5: 'BNB' -> 'BNB'
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '243.97'
$ws.Cells.Item(5, 5).Value = '  +0.71%  '

# Row 6: 'XRP' -> 'XRP'
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.6391'
$ws.Cells.Item(6, 5).Value = '  +3.17%  '

# Row 7: 'USDC' -> 'USDC'
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.000'
$ws.Cells.Item(7, 5).Value = '  -0.07%  '

# Row 8: 'Dogecoin' -> 'OKB'
$ws.Cells.Item(8, 2).Value = 'OKB'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '46.73'
$ws.Cells.Item(8, 5).Value = '  +2.96%  '

# Row 9: 'Cardano' -> 'Dogecoin'
$ws.Cells.Item(9, 2).Value = 'Dogecoin'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07473'
$ws.Cells.Item(9, 5).Value = '  +1.45%  '

# Row 10: 'Solana' -> 'Cardano'
$ws.Cells.Item(10, 2).Value = 'Cardano'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.2979'
$ws.Cells.Item(10, 5).Value = '  +2.25%  '

# Row 11: 'TRON' -> 'Solana'
$ws.Cells.Item(11, 2).Value = 'Solana'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '24.29'
$ws.Cells.Item(11, 5).Value = '  +5.18%  '

# Row 12: 'WrappedEther' -> 'TRON'
$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.07642'
$ws.Cells.Item(12, 5).Value = '  -0.52%  '

# Row 13: 'Polkadot' -> 'WrappedEther'
$ws.Cells.Item(13, 2).Value = 'WrappedEther'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(13, 4).Value = '1.854.23'
$ws.Cells.Item(13, 5).Value = '  +1.74%  '

# Row 14: 'Polygon' -> 'Polkadot'
$ws.Cells.Item(14, 2).Value = 'Polkadot'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.037'
$ws.Cells.Item(14, 5).Value = '  +1.75%  '

# Row 15: 'Litecoin' -> 'Polygon'
$ws.Cells.Item(15, 2).Value = 'Polygon'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.6870'
$ws.Cells.Item(15, 5).Value = '  +3.32%  '

# Row 16: 'ShibaInu' -> 'Litecoin'
$ws.Cells.Item(16, 2).Value = 'Litecoin'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '83.65'
$ws.Cells.Item(16, 5).Value = '  +1.75%  '

# Row 17: 'Uniswap' -> 'ShibaInu'
$ws.Cells.Item(17, 2).Value = 'ShibaInu'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.000009498'
$ws.Cells.Item(17, 5).Value = '  +6.30%  '

# Row 18: 'WrappedBTC' -> 'Uniswap'
$ws.Cells.Item(18, 2).Value = 'Uniswap'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '6.047'
$ws.Cells.Item(18, 5).Value = '  +3.28%  '

# Row 19: 'WrappedliquidstakedEther2.0' -> 'WrappedBTC'
$ws.Cells.Item(19, 2).Value = 'WrappedBTC'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(19, 4).Value = '29.730.61'
$ws.Cells.Item(19, 5).Value = '  +2.00%  '

# Row 20: 'BitcoinCash' -> 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(20, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(20, 4).Value = '2.115.60'
$ws.Cells.Item(20, 5).Value = '  +2.32%  '

# Row 21: 'Avalanche' -> 'BitcoinCash'
$ws.Cells.Item(21, 2).Value = 'BitcoinCash'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '235.18'
$ws.Cells.Item(21, 5).Value = '  -0.73%  '

# Row 22: 'Dai' -> 'Avalanche'
$ws.Cells.Item(22, 2).Value = 'Avalanche'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '12.61'
$ws.Cells.Item(22, 5).Value = '  +1.38%  '

# Row 23: 'Chainlink' -> 'Dai'
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.9999'
$ws.Cells.Item(23, 5).Value = '  -0.06%  '

# Row 24: 'BinanceUSD' -> 'Chainlink'
$ws.Cells.Item(24, 2).Value = 'Chainlink'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '7.387'
$ws.Cells.Item(24, 5).Value = '  +0.65%  '

# Row 25: 'Monero' -> 'BinanceUSD'
$ws.Cells.Item(25, 2).Value = 'BinanceUSD'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '1.001'
$ws.Cells.Item(25, 5).Value = '  -0.01%  '

# Row 26: 'Stellar' -> 'Monero'
$ws.Cells.Item(26, 2).Value = 'Monero'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '158.09'
$ws.Cells.Item(26, 5).Value = '  +0.02%  '

# Row 27: 'Cosmos' -> 'Stellar'
$ws.Cells.Item(27, 2).Value = 'Stellar'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.1416'
$ws.Cells.Item(27, 5).Value = '  +0.48%  '

# Row 28: 'EthereumClassic' -> 'Cosmos'
$ws.Cells.Item(28, 2).Value = 'Cosmos'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '8.476'
$ws.Cells.Item(28, 5).Value = '  -0.07%  '

# Row 29: 'Hedera' -> 'EthereumClassic'
$ws.Cells.Item(29, 2).Value = 'EthereumClassic'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '17.89'
$ws.Cells.Item(29, 5).Value = '  +1.48%  '

# Row 30: 'PancakeSwap' -> 'Hedera'
$ws.Cells.Item(30, 2).Value = 'Hedera'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.06249'
$ws.Cells.Item(30, 5).Value = '  +5.61%  '

# Row 31: 'Toncoin' -> 'PancakeSwap'
$ws.Cells.Item(31, 2).Value = 'PancakeSwap'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.491'
$ws.Cells.Item(31, 5).Value = '  +0.46%  '

# Row 32: 'Filecoin' -> 'Toncoin'
$ws.Cells.Item(32, 2).Value = 'Toncoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.273'
$ws.Cells.Item(32, 5).Value = '  +5.77%  '

# Row 33: 'InternetComputer(DFINITY)' -> 'Filecoin'
$ws.Cells.Item(33, 2).Value = 'Filecoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.142'
$ws.Cells.Item(33, 5).Value = '  +1.46%  '

# Row 34: 'ARBITRUM' -> 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.086'
$ws.Cells.Item(34, 5).Value = '  +0.34%  '

# Row 35: 'LidoDAOToken' -> 'ARBITRUM'
$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.174'
$ws.Cells.Item(35, 5).Value = '  +3.15%  '

# Row 36: 'ImmutableX' -> 'LidoDAOToken'
$ws.Cells.Item(36, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '1.852'
$ws.Cells.Item(36, 5).Value = '  -0.15%  '

# Row 37: 'HuobiToken' -> 'ImmutableX'
$ws.Cells.Item(37, 2).Value = 'ImmutableX'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.7275'
$ws.Cells.Item(37, 5).Value = '  +0.57%  '

# Row 38: 'MXToken' -> 'HuobiToken'
$ws.Cells.Item(38, 2).Value = 'HuobiToken'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '2.605'
$ws.Cells.Item(38, 5).Value = '  +0.15%  '

# Row 39: 'VeChain' -> 'MXToken'
$ws.Cells.Item(39, 2).Value = 'MXToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.845'
$ws.Cells.Item(39, 5).Value = '  -0.08%  '

# Row 40: 'Maker' -> 'VeChain'
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.01778'
$ws.Cells.Item(40, 5).Value = '  +1.81%  '

# Row 41: 'TrustWalletToken' -> 'Maker'
$ws.Cells.Item(41, 2).Value = 'Maker'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(41, 4).Value = '1.200.77'
$ws.Cells.Item(41, 5).Value = '  -1.58%  '

# Row 42: 'FraxShare' -> 'TrustWalletToken'
$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.9237'
$ws.Cells.Item(42, 5).Value = '  +0.47%  '

# Row 43: 'PaxDollar' -> 'FraxShare'
$ws.Cells.Item(43, 2).Value = 'FraxShare'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '6.143'
$ws.Cells.Item(43, 5).Value = '  -2.08%  '

# Row 44: 'RocketPoolETH' -> 'PaxDollar'
$ws.Cells.Item(44, 2).Value = 'PaxDollar'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.9999'
$ws.Cells.Item(44, 5).Value = '  -0.07%  '

# Row 45: 'Quant' -> 'RocketPoolETH'
$ws.Cells.Item(45, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(45, 4).Value = '2.021.15'
$ws.Cells.Item(45, 5).Value = '  +2.58%  '

# Row 46: 'Aave' -> 'Quant'
$ws.Cells.Item(46, 2).Value = 'Quant'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '101.84'
$ws.Cells.Item(46, 5).Value = '  +0.02%  '

# Row 47: 'BabyDogeCoin' -> 'Aave'
$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '65.88'
$ws.Cells.Item(47, 5).Value = '  +1.71%  '

# Row 48: 'TheSandbox' -> 'TheSandbox'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.4052'
$ws.Cells.Item(48, 5).Value = '  +0.88%  '

# Row 49: 'EnergySwap' -> 'EnergySwap'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '9.162'
$ws.Cells.Item(49, 5).Value = '  +0.41%  '

# Row 50: 'Cronos' -> 'Cronos'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.05793'
$ws.Cells.Item(50, 5).Value = '  +0.73%  '

# Row 51: 'RenderToken' -> 'RenderToken'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.642'
$ws.Cells.Item(51, 5).Value = '  +4.23%  '
